# Adds a new "transcripcionFonetica" column (F) to the familias sheet and
# replaces the sample placeholder rows with real family/subfamily data,
# including a row where the "agrupaciones_ids" link field is a single
# numeric id instead of a comma separated list (condición en caso de no
# tener campo a vincular).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column F: "transcripcionFonetica" -----------------------
# This shifts comunidades..normaDeEscritura (old F..U) one column to the
# right (new G..V) and carries the header style (wrap text) along.
$ws.Columns.Item(6).Insert()

$ws.Range("F1").Value = "transcripcionFonetica"
$ws.Columns.Item(6).ColumnWidth = 19.4

# Header row now wraps onto two lines because of the new, longer header.
$ws.Rows.Item(1).RowHeight = 28.35

# --- Row 2: Otomangue -------------------------------------------------
$ws.Range("B2").Value = "Otomangue"
$ws.Range("C2").Value = "Otomangue"
$ws.Range("E2").Value = "1,2,3,4,5,6"

# --- Row 3: Yuto-Nahua / Yuto-azteca / Proto-Nahua ---------------------
$ws.Range("B3").Value = "Yuto-Nahua"
$ws.Range("C3").Value = "Yuto-azteca"
$ws.Range("D3").Value = "Proto-Nahua"
$ws.Range("E3").Value = "7,8,9,10,11,12,13,14,15"

# --- Row 4: Álgica, linked by a single numeric id (no list) -----------
$ws.Range("B4").Value = "Álgica"
$ws.Range("C4").Value = "Álgica"
$ws.Range("E4").Value = 16

# --- Restore cursor / selection ---------------------------------------
$ws.Range("E10").Select()
